$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting the existing row 12 (and all
# following rows) down by one. This matches the dimension growing from
# A1:R90 to A1:R91.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new record.
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(12, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(12, 4).Value = 44847
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = 100112040
$ws.Cells.Item(12, 7).Value = "Cilantro"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 200
$ws.Cells.Item(12, 11).Value = 800
$ws.Cells.Item(12, 12).Value = 1000
$ws.Cells.Item(12, 13).Value = 900
$ws.Cells.Item(12, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(12, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(12, 16).Value = 450
$ws.Cells.Item(12, 17).Value = 2
$ws.Cells.Item(12, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest
# of column D.
$ws.Cells.Item(12, 4).NumberFormat = $ws.Cells.Item(13, 4).NumberFormat
